$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.262.85'
$ws.Range("E2").Value = '  +1.88%  '

$ws.Range("D3").Value = '1.843.90'
$ws.Range("E3").Value = '  +1.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.03'
$ws.Range("E5").Value = '  +1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  +2.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.84'
$ws.Range("E8").Value = '  +14.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.305'
$ws.Range("E9").Value = '  +4.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0694'
$ws.Range("E10").Value = '  +1.60%  '

$ws.Range("E11").Value = '  +3.44%  '

$ws.Range("E12").Value = '  +2.03%  '

$ws.Range("E13").Value = '  +2.11%  '

$ws.Range("D14").Value = '1.839.58'
$ws.Range("E14").Value = '  +1.45%  '

$ws.Range("E15").Value = '  +6.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.659'
$ws.Range("E16").Value = '  +3.70%  '

$ws.Range("D17").Value = '35.282.66'
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.56'
$ws.Range("E18").Value = '  +1.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.00'
$ws.Range("E19").Value = '  +1.04%  '

$ws.Range("D20").Value = '0.0₃0793'
$ws.Range("E20").Value = '  +2.10%  '

$ws.Range("E21").Value = '  +7.62%  '

$ws.Range("E22").Value = '  +13.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").Value = '  -1.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.79'
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.89'
$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.88'
$ws.Range("E27").Value = '  +3.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.122'
$ws.Range("E28").Value = '  +1.07%  '

$ws.Range("B29").Value = 'EURNeutrino'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
$ws.Range("D29").Value = '3.650.26'
$ws.Range("E29").Value = '  +50.24%  '

$ws.Range("B30").Value = 'BinanceUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.33'
$ws.Range("E31").Value = '  +7.73%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.95'
$ws.Range("E32").Value = '  +3.24%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.05'
$ws.Range("E33").Value = '  +2.73%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0537'
$ws.Range("E34").Value = '  +3.52%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.89'
$ws.Range("E35").Value = '  +4.07%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.672'
$ws.Range("E36").Value = '  +2.38%  '

$ws.Range("B37").Value = 'Aave'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '90.11'
$ws.Range("E37").Value = '  +11.07%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.08'
$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.337.98'
$ws.Range("E39").Value = '  -2.12%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  +8.61%  '

$ws.Range("E41").Value = '  +2.83%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("E42").Value = '  +1.34%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.81'
$ws.Range("E43").Value = '  +6.84%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.25'
$ws.Range("E44").Value = '  +6.56%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  +1.53%  '

$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").Value = '  +0.74%  '

$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0518'
$ws.Range("E47").Value = '  +3.33%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.06'
$ws.Range("E48").Value = '  +4.03%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.012.23'
$ws.Range("E49").Value = '  +2.05%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.38'
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.997'
$ws.Range("E51").Value = '  -0.22%  '
